# This workbook's data rows (2-22) are being reshuffled: the "Fecha" (D),
# "Volumen" (M), "Precio mínimo/máximo/promedio ponderado" (N/O/P) and
# "Precio $/Kg" (S) columns move together as a record, while all other
# columns (E,F,G,H,I,J,K,L,Q,R,T) keep their original values per row.
#
# Build the permutation: new row -> source row (using original data).
$map = @{
    2  = 13
    3  = 9
    4  = 20
    5  = 12
    6  = 10
    7  = 17
    8  = 5
    9  = 7
    10 = 16
    11 = 19
    12 = 3
    13 = 2
    14 = 8
    15 = 6
    16 = 15
    17 = 4
    18 = 11
    19 = 14
    20 = 21
    21 = 22
    22 = 18
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original values for the columns that move (D, M, N, O, P, S)
# before any of them get overwritten, since source rows are also targets.
# NOTE: use .Value2 (not .Value) - this runtime mishandles chained/stored
# .Value reads, returning a bogus reflection object instead of the data.
$orig = @{}
foreach ($r in 2..22) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2   # D: Fecha
        M = $ws.Cells.Item($r, 13).Value2  # M: Volumen
        N = $ws.Cells.Item($r, 14).Value2  # N: Precio minimo
        O = $ws.Cells.Item($r, 15).Value2  # O: Precio maximo
        P = $ws.Cells.Item($r, 16).Value2  # P: Precio promedio ponderado
        S = $ws.Cells.Item($r, 19).Value2  # S: Precio $/Kg
    }
}

# Apply new values based on the permutation mapping.
foreach ($r in 2..22) {
    $src = $map[$r]
    $ws.Cells.Item($r, 4).Value2 = $orig[$src].D
    $ws.Cells.Item($r, 13).Value2 = $orig[$src].M
    $ws.Cells.Item($r, 14).Value2 = $orig[$src].N
    $ws.Cells.Item($r, 15).Value2 = $orig[$src].O
    $ws.Cells.Item($r, 16).Value2 = $orig[$src].P
    $ws.Cells.Item($r, 19).Value2 = $orig[$src].S
}
